$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 219-220, pushing the old rows 219:361 down to 221:363
$ws.Rows("219:220").Insert()

# Row 219 (new) - "Primera" quality entry
$ws.Range("A219").Value = 11
$ws.Range("B219").Value = "Vega Monumental Concepción"
$ws.Range("C219").Value = "Bíobío"
$ws.Range("D219").Value = 44960
$ws.Range("E219").Value = 8
$ws.Range("F219").Value = 100112009
$ws.Range("G219").Value = "Acelga"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 100
$ws.Range("K219").Value = 700
$ws.Range("L219").Value = 800
$ws.Range("M219").Value = 750
$ws.Range("N219").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O219").Value = "Región de Ñuble"
$ws.Range("P219").Value = 750
$ws.Range("Q219").Value = 1
$ws.Range("R219").Value = "Hortaliza"

# Row 220 (new) - "Segunda" quality entry
$ws.Range("A220").Value = 11
$ws.Range("B220").Value = "Vega Monumental Concepción"
$ws.Range("C220").Value = "Bíobío"
$ws.Range("D220").Value = 44960
$ws.Range("E220").Value = 8
$ws.Range("F220").Value = 100112009
$ws.Range("G220").Value = "Acelga"
$ws.Range("H220").Value = "Sin especificar"
$ws.Range("I220").Value = "Segunda"
$ws.Range("J220").Value = 50
$ws.Range("K220").Value = 600
$ws.Range("L220").Value = 600
$ws.Range("M220").Value = 600
$ws.Range("N220").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O220").Value = "Región de Ñuble"
$ws.Range("P220").Value = 600
$ws.Range("Q220").Value = 1
$ws.Range("R220").Value = "Hortaliza"
